$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 87: Pimiento - Zafiro rojo - Primera
$ws.Range('A87').Value = 12
$ws.Range('B87').Value = 'Mapocho Venta Directa de Santiago'
$ws.Range('C87').Value = 'Metropolitana'
$ws.Range('D87').Value = 44585
$ws.Range('D87').NumberFormat = $ws.Range('D86').NumberFormat
$ws.Range('E87').Value = 13
$ws.Range('F87').Value = 100112002
$ws.Range('G87').Value = 'Pimiento'
$ws.Range('H87').Value = 'Zafiro rojo'
$ws.Range('I87').Value = 'Primera'
$ws.Range('J87').Value = 25
$ws.Range('K87').Value = 17000
$ws.Range('L87').Value = 17000
$ws.Range('M87').Value = 17000
$ws.Range('N87').Value = '$/caja 18 kilos'
$ws.Range('O87').Value = 'Provincia de Limarí'
$ws.Range('P87').Value = 944
$ws.Range('Q87').Value = 18
$ws.Range('R87').Value = 'Hortaliza'

# Row 88: Pimiento - Zafiro rojo - Segunda
$ws.Range('A88').Value = 12
$ws.Range('B88').Value = 'Mapocho Venta Directa de Santiago'
$ws.Range('C88').Value = 'Metropolitana'
$ws.Range('D88').Value = 44585
$ws.Range('D88').NumberFormat = $ws.Range('D86').NumberFormat
$ws.Range('E88').Value = 13
$ws.Range('F88').Value = 100112002
$ws.Range('G88').Value = 'Pimiento'
$ws.Range('H88').Value = 'Zafiro rojo'
$ws.Range('I88').Value = 'Segunda'
$ws.Range('J88').Value = 30
$ws.Range('K88').Value = 14000
$ws.Range('L88').Value = 14000
$ws.Range('M88').Value = 14000
$ws.Range('N88').Value = '$/caja 18 kilos'
$ws.Range('O88').Value = 'Provincia de Limarí'
$ws.Range('P88').Value = 778
$ws.Range('Q88').Value = 18
$ws.Range('R88').Value = 'Hortaliza'
